$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
}

Replace-Text "2025-03-30 Sunday" "2025-03-31 Monday"

Replace-Text "17×32=544" "23×26=598"
Replace-Text "32×69=2208" "39×84=3276"
Replace-Text "87×60=5220" "21×86=1806"
Replace-Text "14×99=1386" "79×89=7031"
Replace-Text "78×59=4602" "46×78=3588"

Replace-Text "83×92=7636" "31×94=2914"
Replace-Text "18×42=756" "69×86=5934"
Replace-Text "17×72=1224" "32×53=1696"
Replace-Text "35×44=1540" "55×94=5170"
Replace-Text "81×39=3159" "90×82=7380"

Replace-Text "67×86=5762" "55×70=3850"
Replace-Text "14×49=686" "46×50=2300"
Replace-Text "42×28=1176" "54×20=1080"
Replace-Text "38×14=532" "98×21=2058"
Replace-Text "76×83=6308" "14×16=224"

Replace-Text "60×49=2940" "70×33=2310"
Replace-Text "39×74=2886" "45×95=4275"
Replace-Text "71×90=6390" "41×89=3649"
Replace-Text "71×12=852" "22×82=1804"
Replace-Text "31×23=713" "18×38=684"

Replace-Text "55×73=4015" "46×94=4324"
Replace-Text "33×13=429" "87×87=7569"
Replace-Text "13×85=1105" "34×92=3128"
Replace-Text "56×85=4760" "93×24=2232"
Replace-Text "95×88=8360" "12×41=492"
